$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: Cost (D1), Total (F1) and the grand-total formula (G1) ---
$ws.Range("D1").Value = "Cost"
$ws.Range("F1").Value = "Total"
$ws.Range("G1").Formula = "=SUMPRODUCT(C2:C9,D2:D9)"

# --- Per-line "Cost" values (column D) ---
$ws.Range("D2").Value = 0.1315
$ws.Range("D3").Value = 0.0143
$ws.Range("D4").Value = 0.0149
$ws.Range("D5").Value = 0.0146
$ws.Range("D6").Value = 0.0018
$ws.Range("D7").Value = 0.0017
$ws.Range("D8").Value = 0.0016
$ws.Range("D9").Value = 0.0152

# --- Row 7: JLCPCB ref for the 100kΩ Resistor changed ---
$ws.Range("B7").Value = "C149504"

# --- Row 8 is now the 3,3kΩ Resistor (replacing the old 3.3kΩ Resistor line) ---
$ws.Range("A8").Value = "3,3kΩ Resistor"
$ws.Range("B8").Value = "C26010"

# --- Apply the "Comma" cell style + custom currency number format to the total ---
$ws.Range("G1").Style = "Comma"
$ws.Range("G1").NumberFormat = "[$$-409]#,##0.0000"

# --- Column widths for the new Cost / Total columns ---
# (values chosen so the saved width lands as close as possible to the
#  target 18.83203125 / 41.33203125 character widths)
$ws.Columns("E").ColumnWidth = 18
$ws.Columns("G").ColumnWidth = 40.5

# --- Selection moves to the new total cell ---
$ws.Range("G1").Select()
